$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Sheet1") becomes "datos" -------------------------------------
$datos = $wb.Worksheets.Item(1)
$datos.Name = "datos"

# --- Add the new "metadatos" sheet right after "datos" ----------------------
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $datos)
$meta.Name = "metadatos"

# --- Header row --------------------------------------------------------------
$meta.Range("A1").Value = "Variables"
$meta.Range("B1").Value = "Descripción"
$meta.Range("C1").Value = "Fuente"
$meta.Range("D1").Value = "Fecha_de_extracción"

# --- Row 2: anno / Año --------------------------------------------------------
$meta.Range("A2").Value = "anno"
$meta.Range("B2").Value = "Año"
$meta.Range("C2").Value = "…"
$meta.Range("D2").Value = 45722

# --- Row 3: codmpio / Código del municipio -----------------------------------
$meta.Range("A3").Value = "codmpio"
$meta.Range("B3").Value = "Código del municipio"
$meta.Range("C3").Value = "…"
$meta.Range("D3").Value = 45722

# --- Row 4: secu / descripción larga / fuente --------------------------------
$meta.Range("A4").Value = "secu"
$meta.Range("B4").Value = "Número de menores de edad secuestrados en el contexto del conflicto"
$meta.Range("C4").Value = "Panel CEDE - Registro Único de Víctimas`nInstituto Colombiano de Bienestar Familiar (ICBF)"
$meta.Range("D4").Value = 45722

# --- Formatting ---------------------------------------------------------------
# Regular (non-default-scheme) Calibri font used across the table
$meta.Range("A1:D4").Font.Name = "Calibri"

# Date column gets Excel's built-in short date format (numFmtId 15)
$meta.Range("D2:D4").NumberFormat = "d-mmm-yy"

# "secu" (A4) carries its own distinct font in the source workbook
$meta.Range("A4").Font.Name = "Arial"
$meta.Range("A4").Font.Name = "Calibri"

# The multi-line source text in C4 auto-grows the row; put it back to the
# sheet's normal (default) height so no explicit row height is stored.
$meta.Rows.Item(4).AutoFit()

# --- Selections / active views ------------------------------------------------
$datos.Range("C1").Select()
$meta.Range("G8").Select()
